# The edit adds a new "OD = 1" biomass-conversion helper column (I) on the
# "Trial 1" sheet: a label, a density-ratio formula, and a units label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trial 1")

$ws.Range("I8").Value = "OD = 1"
$ws.Range("I9").Formula = "=0.69/0.895"
$ws.Range("I10").Value = "mg/ml"

$ws.Range("I11").Select()
